$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$people = @(
    @{ Name = "Dennis"; Email = "den@nis.nis" },
    @{ Name = "Cooker"; Email = "cook@er.errr" },
    @{ Name = "Test";   Email = "tester@test.se" },
    @{ Name = "Jack";   Email = "jacklo@ers.se" }
)

$row = 5
foreach ($p in $people) {
    $ws.Range("A$row").Value = $p.Name
    $ws.Range("B$row").Value = $p.Email
    $ws.Range("C$row").Value = "nomail"
    $ws.Hyperlinks.Add($ws.Range("B$row"), "mailto:" + $p.Email)
    $row++
}

$ws.Range("A9").Select()
